# Update the "F" column (numeric counts) for matching rows on both the
# "展览" and "全部类型" sheets, which contain duplicated data for rows 2-20.

$wb = $excel.ActiveWorkbook

# Row (by sheet row number) -> new value for column F
$updates = @{
    2  = 117
    4  = 11836
    5  = 1029
    6  = 122
    9  = 154
    10 = 173
    11 = 29
    12 = 52
    13 = 55
    15 = 35
    16 = 347
    17 = 1422
    18 = 80
    19 = 916
    20 = 115
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
